$d = $word.ActiveDocument

# The document uses "{{Placeholder}}" style tokens that need to become
# single-brace "{Placeholder}" tokens. Literal Find text that mixes an
# opening "{" and closing "}" can't be matched reliably with a plain
# (non-wildcard) search, so we use wildcard character classes ("[{]"/"[}]")
# to match the braces literally while keeping the rest of the text literal.

# "HEI {{Name}}" -> "HEI {Name}"
# (The run carries w:caps formatting, so the text must be matched in its
# rendered, upper-cased form.)
$r1 = $d.Content.Find.Execute(
    "HEI [{][{]NAME[}][}]",
    $true, $false, $true, $false, $false, $true, 1, $false,
    "HEI {Name}", 2)

# "{{ActivationCode}}" -> "{ActivationCode}"
$r2 = $d.Content.Find.Execute(
    "[{][{]ActivationCode[}][}]",
    $true, $false, $true, $false, $false, $true, 1, $false,
    "{ActivationCode}", 2)

Write-Output "Name token replaced: $r1"
Write-Output "ActivationCode token replaced: $r2"
